$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 140
$ws.Range("I9").Value = 148.57143
$ws.Range("J9").Value = 110
$ws.Range("K9").Value = 148.57143
$ws.Range("L9").Value = 110
$ws.Range("M9").Value = 20.42857000000001
$ws.Range("N9").Value = -448
$ws.Range("H38").Value = 654.0769
$ws.Range("I38").Value = 217.44444
$ws.Range("J38").Value = 1636.5
$ws.Range("K38").Value = 652.33332
$ws.Range("L38").Value = 4909.5
$ws.Range("M38").Value = -280.33332
$ws.Range("N38").Value = -5653.5
$ws.Range("H53").Value = 207.53334
$ws.Range("I53").Value = 201.22223
$ws.Range("J53").Value = 217
$ws.Range("K53").Value = 201.22223
$ws.Range("L53").Value = 217
$ws.Range("M53").Value = 435.77777
$ws.Range("N53").Value = -1491
$ws.Range("H58").Value = 2056.647
$ws.Range("I58").Value = 1523.7778
$ws.Range("J58").Value = 2656.125
$ws.Range("K58").Value = 4571.3334
$ws.Range("L58").Value = 7968.375
$ws.Range("M58").Value = -4421.3334
$ws.Range("N58").Value = -8268.375
$ws.Range("H87").Value = 24089.334
$ws.Range("J87").Value = 24089.334
$ws.Range("L87").Value = 24089.334
$ws.Range("N87").Value = -26585.334
$ws.Range("H90").Value = 24089.334
$ws.Range("J90").Value = 24089.334
$ws.Range("L90").Value = 72268.00199999999
$ws.Range("N90").Value = -84748.00199999999
$ws.Range("H99").Value = 4617.8
$ws.Range("I99").Value = 704
$ws.Range("J99").Value = 7227
$ws.Range("K99").Value = 2112
$ws.Range("L99").Value = 21681
$ws.Range("M99").Value = -614
$ws.Range("N99").Value = -24677
$ws.Range("H103").Value = 56413.445
$ws.Range("I103").Value = 100344.9
$ws.Range("J103").Value = 1499.125
$ws.Range("K103").Value = 301034.7
$ws.Range("L103").Value = 4497.375
$ws.Range("M103").Value = -300448.7
$ws.Range("N103").Value = -5669.375
$ws.Range("H138").Value = 2863.48
$ws.Range("I138").Value = 1075.8125
$ws.Range("J138").Value = 3203.988
$ws.Range("K138").Value = 3227.4375
$ws.Range("L138").Value = 9611.964
$ws.Range("M138").Value = 1912.5625
$ws.Range("N138").Value = -19891.964

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 638.7857
$ws.Range("I2").Value = 652
$ws.Range("J2").Value = 605.75
$ws.Range("K2").Value = 652
$ws.Range("L2").Value = 605.75
$ws.Range("M2").Value = -539
$ws.Range("N2").Value = -831.75
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 12824051
$ws.Range("I61").Value = 19609978
$ws.Range("J61").Value = 6188.8887
$ws.Range("K61").Value = 19609978
$ws.Range("L61").Value = 6188.8887
$ws.Range("M61").Value = -19609766
$ws.Range("N61").Value = -6612.8887
$ws.Range("H116").Value = 638.7857
$ws.Range("I116").Value = 652
$ws.Range("J116").Value = 605.75
$ws.Range("K116").Value = 652
$ws.Range("L116").Value = 605.75
$ws.Range("M116").Value = 1642
$ws.Range("N116").Value = -5193.75
$ws.Range("H133").Value = 55464.4
$ws.Range("J133").Value = 55464.4
$ws.Range("L133").Value = 55464.4
$ws.Range("N133").Value = -60524.4
$ws.Range("H136").Value = 12824051
$ws.Range("I136").Value = 19609978
$ws.Range("J136").Value = 6188.8887
$ws.Range("K136").Value = 58829934
$ws.Range("L136").Value = 18566.6661
$ws.Range("M136").Value = -58827384
$ws.Range("N136").Value = -23666.6661

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 638.7857
$ws.Range("I3").Value = 652
$ws.Range("J3").Value = 605.75
$ws.Range("K3").Value = 652
$ws.Range("L3").Value = 605.75
$ws.Range("M3").Value = -538
$ws.Range("N3").Value = -833.75
$ws.Range("H94").Value = 247.18182
$ws.Range("I94").Value = 240.9
$ws.Range("J94").Value = 310
$ws.Range("K94").Value = 240.9
$ws.Range("L94").Value = 310
$ws.Range("M94").Value = 210.1
$ws.Range("N94").Value = -1212
$ws.Range("H107").Value = 72643.5
$ws.Range("I107").Value = 101091
$ws.Range("J107").Value = 1524.75
$ws.Range("K107").Value = 101091
$ws.Range("L107").Value = 1524.75
$ws.Range("M107").Value = -99171
$ws.Range("N107").Value = -5364.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3677316.2
$ws.Range("I107").Value = 6250489.5
$ws.Range("J107").Value = 1354.2858
$ws.Range("K107").Value = 6250489.5
$ws.Range("L107").Value = 1354.2858
$ws.Range("M107").Value = -6248569.5
$ws.Range("N107").Value = -5194.2858

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3952.2222
$ws.Range("I63").Value = 2714
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 8142
$ws.Range("L63").Value = 16500
$ws.Range("M63").Value = -7393
$ws.Range("N63").Value = -17998
$ws.Range("H66").Value = 3952.2222
$ws.Range("I66").Value = 2714
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 24426
$ws.Range("L66").Value = 49500
$ws.Range("M66").Value = -20682
$ws.Range("N66").Value = -56988
$ws.Range("H131").Value = 3476.2954
$ws.Range("I131").Value = 405
$ws.Range("J131").Value = 5231.3213
$ws.Range("K131").Value = 1215
$ws.Range("L131").Value = 15693.9639
$ws.Range("M131").Value = 3825
$ws.Range("N131").Value = -25773.9639
$ws.Range("H132").Value = 2730.075
$ws.Range("I132").Value = 2426.2
$ws.Range("J132").Value = 3033.95
$ws.Range("K132").Value = 21835.8
$ws.Range("L132").Value = 27305.55
$ws.Range("M132").Value = -19305.8
$ws.Range("N132").Value = -32365.55
$ws.Range("H137").Value = 40684.93
$ws.Range("I137").Value = 8549.8125
$ws.Range("J137").Value = 80235.84
$ws.Range("K137").Value = 25649.4375
$ws.Range("L137").Value = 240707.52
$ws.Range("M137").Value = -20549.4375
$ws.Range("N137").Value = -250907.52

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 341.11765
$ws.Range("I107").Value = 339.93332
$ws.Range("K107").Value = 339.93332
$ws.Range("M107").Value = 1580.06668

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 8647.875
$ws.Range("J93").Value = 3234.6667
$ws.Range("L93").Value = 3234.6667
$ws.Range("N93").Value = -5730.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 645.8
$ws.Range("I107").Value = 580.25
$ws.Range("J107").Value = 720.7143
$ws.Range("K107").Value = 1740.75
$ws.Range("L107").Value = 2162.1429
$ws.Range("M107").Value = 179.25
$ws.Range("N107").Value = -6002.1429
$ws.Range("H109").Value = 2692933.2
$ws.Range("J109").Value = 39400
$ws.Range("L109").Value = 39400
$ws.Range("N109").Value = -42174
$ws.Range("H138").Value = 75723.25
$ws.Range("J138").Value = 75723.25
$ws.Range("L138").Value = 75723.25
$ws.Range("N138").Value = -86003.25
